$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Thresholds by uom")

# --- Header renames (A1/B1): "Simple_Parameter"/"UOM" -> "Simple Parameter"/"uom" ---
$ws.Range("A1").Value = "Simple Parameter"
$ws.Range("B1").Value = "uom"

# --- Clear the placeholder "NA" values in the unused Fresh_2 / Marine_1 / Marine_2
#     columns, leaving the cell formatting (style) intact but with no content. ---
$naCells = @(
  "E2","F2","G2","H2","I2","J2",
  "E3","F3","G3","H3","I3","J3",
  "E4","F4","G4","H4","I4","J4",
  "E5","F5","G5","H5","I5","J5",
  "E6","F6","I6","J6",
  "E7","F7","I7","J7",
  "E11","F11",
  "E12","F12",
  "E13","F13",
  "G14","H14","I14","J14",
  "G15","H15","I15","J15",
  "G16","H16","I16","J16",
  "G17","H17","I17","J17",
  "G18","H18","I18","J18",
  "G19","H19","I19","J19",
  "G20","H20","I20","J20",
  "C21","D21","E21","F21","I21","J21",
  "C22","D22","E22","F22","I22","J22",
  "C23","D23","E23","F23","I23","J23",
  "C24","D24","E24","F24","I24","J24",
  "C25","D25","E25","F25","I25","J25",
  "C26","D26","E26","F26","I26","J26",
  "C27","D27","E27","F27","I27","J27",
  "C28","D28","E28","F28","I28","J28",
  "C29","D29","E29","F29","I29","J29"
)

foreach ($addr in $naCells) {
  $ws.Range($addr).Value = $null
}

# --- Selection / active cell on the sheet moves from F16 to D21:D29 ---
$ws.Activate() | Out-Null
$ws.Range("D21:D29").Select() | Out-Null
